$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "status_code: 406"
$ws.Range("C4").Value = "status_code: 406"

$ws.Range("I5").Value = "HTTPConnectionPool(host='10.1.2.249', port=80): Max retries exceeded with url: / (Caused by ConnectTimeoutError(<urllib3.connection.HTTPConnection object at 0x000001A293AF8EE0>, 'Connection to 10.1.2.249 timed out. (connect timeout=10)'))"
